$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarters: Dec-2018, Sep-2018),
# shifting the existing quarterly data from D:K to F:M.
$ws.Columns("D:E").Insert()

# Carry over number/date formatting from column F/G (the old column D/E) to
# the two newly inserted blank columns so the new cells keep the right style.
$ws.Range("F5:G102").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Populate the two new columns (D = 31-Dec-2018, E = 30-Sep-2018) with
# the new quarter's figures for every existing line item.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 679000
$ws.Range("E8").Value = 655000
$ws.Range("D9").Value = 986900
$ws.Range("E9").Value = 459200
$ws.Range("D10").Value = -307900
$ws.Range("E10").Value = 195800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 14900
$ws.Range("E15").Value = 17100
$ws.Range("D17").Value = 577200
$ws.Range("E17").Value = 557100
$ws.Range("D18").Value = 101800
$ws.Range("E18").Value = 97900
$ws.Range("D20").Value = 1600
$ws.Range("E20").Value = 1300
$ws.Range("D21").Value = 118300
$ws.Range("E21").Value = 116300
$ws.Range("D22").Value = 3200
$ws.Range("E22").Value = 3200
$ws.Range("D23").Value = 100200
$ws.Range("E23").Value = 96100
$ws.Range("D24").Value = 12000
$ws.Range("E24").Value = 11500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 88200
$ws.Range("E26").Value = 84500
$ws.Range("D27").Value = 88200
$ws.Range("E27").Value = 84500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1600
$ws.Range("E32").Value = -1300
$ws.Range("D33").Value = 88200
$ws.Range("E33").Value = 84500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 88200
$ws.Range("E35").Value = 84500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 395900
$ws.Range("E41").Value = 426200
$ws.Range("D42").Value = 59900
$ws.Range("E42").Value = 65300
$ws.Range("D43").Value = 837600
$ws.Range("E43").Value = 771700
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 36800
$ws.Range("E45").Value = 42300
$ws.Range("D46").Value = 1330200
$ws.Range("E46").Value = 1305500
$ws.Range("D47").Value = 27000
$ws.Range("E47").Value = 19000
$ws.Range("D48").Value = 158700
$ws.Range("E48").Value = 154500
$ws.Range("D49").Value = 810300
$ws.Range("E49").Value = 815000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 28100
$ws.Range("E52").Value = 31600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2354300
$ws.Range("E54").Value = 2325700
$ws.Range("D57").Value = 13300
$ws.Range("E57").Value = 15300
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 597300
$ws.Range("E59").Value = 586800
$ws.Range("D60").Value = 610600
$ws.Range("E60").Value = 602100
$ws.Range("D61").Value = 349300
$ws.Range("E61").Value = 349200
$ws.Range("D62").Value = 40100
$ws.Range("E62").Value = 39900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1000000
$ws.Range("E66").Value = 991200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 889300
$ws.Range("E72").Value = 873100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1354300
$ws.Range("E76").Value = 1334500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 88200
$ws.Range("E81").Value = 84500
$ws.Range("D83").Value = 14900
$ws.Range("E83").Value = 17100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 60900
$ws.Range("E89").Value = 128100
$ws.Range("D91").Value = -20000
$ws.Range("E91").Value = -11200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -20300
$ws.Range("E94").Value = 7700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -70100
$ws.Range("E100").Value = 1800
$ws.Range("D101").Value = -900
$ws.Range("E101").Value = -600
$ws.Range("D102").Value = -30400
$ws.Range("E102").Value = 137000

